$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H100").Value = 16669017
$ws.Range("I100").Value = 25642604
$ws.Range("J100").Value = 3785.8572
$ws.Range("K100").Value = 25642604
$ws.Range("L100").Value = 3785.8572
$ws.Range("M100").Value = -25642063
$ws.Range("N100").Value = -4867.8572
$ws.Range("H112").Value = 1588.7587
$ws.Range("I112").Value = 600
$ws.Range("J112").Value = 1624.0714
$ws.Range("K112").Value = 1800
$ws.Range("L112").Value = 4872.2142
$ws.Range("M112").Value = -692
$ws.Range("N112").Value = -7088.2142
$ws.Range("H116").Value = 2008.5294
$ws.Range("I116").Value = 1730.625
$ws.Range("J116").Value = 2255.5557
$ws.Range("K116").Value = 1730.625
$ws.Range("L116").Value = 2255.5557
$ws.Range("M116").Value = 1711.375
$ws.Range("N116").Value = -9139.555700000001
$ws.Range("H125").Value = 671
$ws.Range("J125").Value = 660
$ws.Range("L125").Value = 5940
$ws.Range("N125").Value = -10860
$ws.Range("H132").Value = 36894.277
$ws.Range("I132").Value = 45523.22
$ws.Range("J132").Value = 3816.6667
$ws.Range("K132").Value = 136569.66
$ws.Range("L132").Value = 11450.0001
$ws.Range("M132").Value = -134039.66
$ws.Range("N132").Value = -16510.0001
$ws.Range("H137").Value = 936.62964
$ws.Range("I137").Value = 707.5333000000001
$ws.Range("J137").Value = 1223
$ws.Range("K137").Value = 2122.5999
$ws.Range("L137").Value = 3669
$ws.Range("M137").Value = 427.4000999999998
$ws.Range("N137").Value = -8769

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 312617.5
$ws.Range("I32").Value = 2936.6626
$ws.Range("J32").Value = 2882968.5
$ws.Range("K32").Value = 2936.6626
$ws.Range("L32").Value = 2882968.5
$ws.Range("M32").Value = -2649.6626
$ws.Range("N32").Value = -2883542.5
$ws.Range("H61").Value = 1136.6666
$ws.Range("I61").Value = 939.2727
$ws.Range("J61").Value = 1446.8572
$ws.Range("K61").Value = 939.2727
$ws.Range("L61").Value = 1446.8572
$ws.Range("M61").Value = -727.2727
$ws.Range("N61").Value = -1870.8572
$ws.Range("H74").Value = 1048.081
$ws.Range("I74").Value = 898.8261
$ws.Range("J74").Value = 1293.2858
$ws.Range("K74").Value = 898.8261
$ws.Range("L74").Value = 1293.2858
$ws.Range("M74").Value = -24.8261
$ws.Range("N74").Value = -3041.2858
$ws.Range("H77").Value = 1048.081
$ws.Range("I77").Value = 898.8261
$ws.Range("J77").Value = 1293.2858
$ws.Range("K77").Value = 4494.1305
$ws.Range("L77").Value = 6466.429
$ws.Range("M77").Value = -126.1305000000002
$ws.Range("N77").Value = -15202.429
$ws.Range("H122").Value = 46730.168
$ws.Range("I122").Value = 26253.84
$ws.Range("J122").Value = 93267.27
$ws.Range("K122").Value = 78761.52
$ws.Range("L122").Value = 279801.81
$ws.Range("M122").Value = -76311.52
$ws.Range("N122").Value = -284701.81
$ws.Range("H136").Value = 1136.6666
$ws.Range("I136").Value = 939.2727
$ws.Range("J136").Value = 1446.8572
$ws.Range("K136").Value = 2817.8181
$ws.Range("L136").Value = 4340.571599999999
$ws.Range("M136").Value = -267.8181
$ws.Range("N136").Value = -9440.571599999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 3108.513
$ws.Range("I20").Value = 2575.3928
$ws.Range("J20").Value = 4465.5454
$ws.Range("K20").Value = 2575.3928
$ws.Range("L20").Value = 4465.5454
$ws.Range("M20").Value = -2328.3928
$ws.Range("N20").Value = -4959.5454
$ws.Range("H52").Value = 24933.334
$ws.Range("I52").Value = 10000
$ws.Range("J52").Value = 26800
$ws.Range("K52").Value = 10000
$ws.Range("L52").Value = 26800
$ws.Range("M52").Value = -9737
$ws.Range("N52").Value = -27326
$ws.Range("H107").Value = 926.96295
$ws.Range("I107").Value = 657.41174
$ws.Range("J107").Value = 1385.2
$ws.Range("K107").Value = 657.41174
$ws.Range("L107").Value = 1385.2
$ws.Range("M107").Value = 1262.58826
$ws.Range("N107").Value = -5225.2
$ws.Range("H121").Value = 24933.334
$ws.Range("I121").Value = 10000
$ws.Range("J121").Value = 26800
$ws.Range("K121").Value = 10000
$ws.Range("L121").Value = 26800
$ws.Range("M121").Value = -8253
$ws.Range("N121").Value = -30294

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 294564.75
$ws.Range("I99").Value = 406142.97
$ws.Range("J99").Value = 1671.875
$ws.Range("K99").Value = 406142.97
$ws.Range("L99").Value = 1671.875
$ws.Range("M99").Value = -404644.97
$ws.Range("N99").Value = -4667.875
$ws.Range("H126").Value = 294564.75
$ws.Range("I126").Value = 406142.97
$ws.Range("J126").Value = 1671.875
$ws.Range("K126").Value = 1218428.91
$ws.Range("L126").Value = 5015.625
$ws.Range("M126").Value = -1215958.91
$ws.Range("N126").Value = -9955.625

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H105").Value = 6676867.5
$ws.Range("J105").Value = 11118604
$ws.Range("L105").Value = 33355812
$ws.Range("N105").Value = -33361054
$ws.Range("H106").Value = 2525
$ws.Range("J106").Value = 2525
$ws.Range("L106").Value = 7575
$ws.Range("N106").Value = -9467
$ws.Range("H113").Value = 728.45
$ws.Range("I113").Value = 584.9
$ws.Range("J113").Value = 764.3375
$ws.Range("K113").Value = 1754.7
$ws.Range("L113").Value = 2293.0125
$ws.Range("M113").Value = 415.3000000000002
$ws.Range("N113").Value = -6633.0125

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 376.53333
$ws.Range("I107").Value = 334.625
$ws.Range("J107").Value = 424.42856
$ws.Range("K107").Value = 334.625
$ws.Range("L107").Value = 424.42856
$ws.Range("M107").Value = 1585.375
$ws.Range("N107").Value = -4264.42856
$ws.Range("H113").Value = 1933.3334
$ws.Range("J113").Value = 2150
$ws.Range("L113").Value = 2150
$ws.Range("N113").Value = -6490

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 2599.5
$ws.Range("I40").Value = 2499.1667
$ws.Range("J40").Value = 2750
$ws.Range("K40").Value = 2499.1667
$ws.Range("L40").Value = 2750
$ws.Range("M40").Value = -2363.1667
$ws.Range("N40").Value = -3022
$ws.Range("H46").Value = 15585.857
$ws.Range("I46").Value = 1400.3334
$ws.Range("J46").Value = 26225
$ws.Range("K46").Value = 1400.3334
$ws.Range("L46").Value = 26225
$ws.Range("M46").Value = -1212.3334
$ws.Range("N46").Value = -26601
$ws.Range("H82").Value = 1787.8
$ws.Range("I82").Value = 2325.7144
$ws.Range("J82").Value = 1317.125
$ws.Range("K82").Value = 2325.7144
$ws.Range("L82").Value = 1317.125
$ws.Range("M82").Value = -1964.7144
$ws.Range("N82").Value = -2039.125
$ws.Range("H85").Value = 1787.8
$ws.Range("I85").Value = 2325.7144
$ws.Range("J85").Value = 1317.125
$ws.Range("K85").Value = 2325.7144
$ws.Range("L85").Value = 1317.125
$ws.Range("M85").Value = -1077.7144
$ws.Range("N85").Value = -3813.125
$ws.Range("H122").Value = 2333.6785
$ws.Range("I122").Value = 2088.7856
$ws.Range("J122").Value = 2578.5715
$ws.Range("K122").Value = 6266.3568
$ws.Range("L122").Value = 7735.7145
$ws.Range("M122").Value = -3816.3568
$ws.Range("N122").Value = -12635.7145
$ws.Range("H136").Value = 13145
$ws.Range("I136").Value = 18248
$ws.Range("J136").Value = 8771
$ws.Range("K136").Value = 54744
$ws.Range("L136").Value = 26313
$ws.Range("M136").Value = -52194
$ws.Range("N136").Value = -31413

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H46").Value = 39804.75
$ws.Range("J46").Value = 39804.75
$ws.Range("L46").Value = 39804.75
$ws.Range("N46").Value = -40266.75
$ws.Range("H122").Value = 1259.4193
$ws.Range("I122").Value = 1229.08
$ws.Range("J122").Value = 1385.8334
$ws.Range("K122").Value = 3687.24
$ws.Range("L122").Value = 4157.5002
$ws.Range("M122").Value = -1237.24
$ws.Range("N122").Value = -9057.5002
$ws.Range("H134").Value = 39804.75
$ws.Range("J134").Value = 39804.75
$ws.Range("L134").Value = 119414.25
$ws.Range("N134").Value = -124484.25
